$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for Application ID 3 (Project ID 4, T2345678D, Pending) entirely.
# This shifts the old row 5 (Application ID 4) up to row 4.
$ws.Rows.Item(4).Delete()

# Update the previously "Deleted" statuses to "Pending"
$ws.Range("D2").Value = "Pending"
$ws.Range("D3").Value = "Pending"

# Update the remaining application (now on row 4) to "Successful" with a new timestamp
$ws.Range("D4").Value = "Successful"
$ws.Range("F4").Value = 45771.755878125

# Update selection to match target state
$ws.Range("D4").Select()
